$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 937.0952
$ws.Range("I92").Value = 696.8182
$ws.Range("J92").Value = 1201.4
$ws.Range("K92").Value = 696.8182
$ws.Range("L92").Value = 1201.4
$ws.Range("M92").Value = 551.1818
$ws.Range("N92").Value = -3697.4
$ws.Range("H100").Value = 2287.375
$ws.Range("I100").Value = 2249.75
$ws.Range("J100").Value = 2325
$ws.Range("K100").Value = 2249.75
$ws.Range("L100").Value = 2325
$ws.Range("M100").Value = -1708.75
$ws.Range("N100").Value = -3407
$ws.Range("H113").Value = 2769.2
$ws.Range("J113").Value = 2742.4
$ws.Range("L113").Value = 2742.4
$ws.Range("N113").Value = -9250.4
$ws.Range("H132").Value = 10426480
$ws.Range("I132").Value = 11912991
$ws.Range("K132").Value = 35738973
$ws.Range("M132").Value = -35736443

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 670.64
$ws.Range("J2").Value = 924.6
$ws.Range("L2").Value = 924.6
$ws.Range("N2").Value = -1150.6
$ws.Range("H116").Value = 670.64
$ws.Range("J116").Value = 924.6
$ws.Range("L116").Value = 924.6
$ws.Range("N116").Value = -5512.6
$ws.Range("H122").Value = 1807.2
$ws.Range("I122").Value = 1858.1364
$ws.Range("J122").Value = 1721
$ws.Range("K122").Value = 5574.4092
$ws.Range("L122").Value = 5163
$ws.Range("M122").Value = -3124.4092
$ws.Range("N122").Value = -10063
$ws.Range("H132").Value = 2450.484
$ws.Range("I132").Value = 1948.9131
$ws.Range("K132").Value = 5846.7393
$ws.Range("M132").Value = -3316.7393

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 670.64
$ws.Range("J3").Value = 924.6
$ws.Range("L3").Value = 924.6
$ws.Range("N3").Value = -1152.6
$ws.Range("H86").Value = 2210.7307
$ws.Range("I86").Value = 2044.0555
$ws.Range("J86").Value = 2585.75
$ws.Range("K86").Value = 2044.0555
$ws.Range("L86").Value = 2585.75
$ws.Range("M86").Value = -921.0554999999999
$ws.Range("N86").Value = -4831.75
$ws.Range("H89").Value = 2210.7307
$ws.Range("I89").Value = 2044.0555
$ws.Range("J89").Value = 2585.75
$ws.Range("K89").Value = 10220.2775
$ws.Range("L89").Value = 12928.75
$ws.Range("M89").Value = -4604.2775
$ws.Range("N89").Value = -24160.75
$ws.Range("H94").Value = 10000256
$ws.Range("I94").Value = 14706096
$ws.Range("J94").Value = 344.875
$ws.Range("K94").Value = 14706096
$ws.Range("L94").Value = 344.875
$ws.Range("M94").Value = -14705645
$ws.Range("N94").Value = -1246.875
$ws.Range("H99").Value = 47620172
$ws.Range("I99").Value = 52632680
$ws.Range("J99").Value = 1350
$ws.Range("K99").Value = 52632680
$ws.Range("L99").Value = 1350
$ws.Range("M99").Value = -52631182
$ws.Range("N99").Value = -4346
$ws.Range("H107").Value = 1254.4117
$ws.Range("I107").Value = 809
$ws.Range("J107").Value = 2071
$ws.Range("K107").Value = 809
$ws.Range("L107").Value = 2071
$ws.Range("M107").Value = 1111
$ws.Range("N107").Value = -5911

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 1629.25
$ws.Range("I17").Value = 505.66666
$ws.Range("J17").Value = 5000
$ws.Range("K17").Value = 505.66666
$ws.Range("L17").Value = 5000
$ws.Range("M17").Value = -331.66666
$ws.Range("N17").Value = -5348
$ws.Range("H31").Value = 1624.875
$ws.Range("I31").Value = 1499.8572
$ws.Range("K31").Value = 1499.8572
$ws.Range("M31").Value = -1204.8572
$ws.Range("H34").Value = 1624.875
$ws.Range("I34").Value = 1499.8572
$ws.Range("K34").Value = 1499.8572
$ws.Range("M34").Value = -1297.8572
$ws.Range("H88").Value = 7000
$ws.Range("J88").Value = 7000
$ws.Range("L88").Value = 7000
$ws.Range("N88").Value = -7812
$ws.Range("H91").Value = 7000
$ws.Range("J91").Value = 7000
$ws.Range("L91").Value = 7000
$ws.Range("N91").Value = -9808
$ws.Range("H99").Value = 1689.5
$ws.Range("I99").Value = 1670.4
$ws.Range("J99").Value = 1721.3334
$ws.Range("K99").Value = 1670.4
$ws.Range("L99").Value = 1721.3334
$ws.Range("M99").Value = -172.4000000000001
$ws.Range("N99").Value = -4717.3334
$ws.Range("H122").Value = 4929.64
$ws.Range("I122").Value = 5051.7085
$ws.Range("K122").Value = 15155.1255
$ws.Range("M122").Value = -12705.1255
$ws.Range("H126").Value = 1689.5
$ws.Range("I126").Value = 1670.4
$ws.Range("J126").Value = 1721.3334
$ws.Range("K126").Value = 5011.200000000001
$ws.Range("L126").Value = 5164.0002
$ws.Range("M126").Value = -2541.200000000001
$ws.Range("N126").Value = -10104.0002
$ws.Range("H132").Value = 1794.2858
$ws.Range("I132").Value = 1114
$ws.Range("J132").Value = 2899.75
$ws.Range("K132").Value = 3342
$ws.Range("L132").Value = 8699.25
$ws.Range("M132").Value = -812
$ws.Range("N132").Value = -13759.25
$ws.Range("H134").Value = 17859158
$ws.Range("I134").Value = 2005.4783
$ws.Range("J134").Value = 100002056
$ws.Range("K134").Value = 6016.4349
$ws.Range("L134").Value = 300006168
$ws.Range("M134").Value = -3481.4349
$ws.Range("N134").Value = -300011238

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 507.43478
$ws.Range("I97").Value = 474.93332
$ws.Range("J97").Value = 568.375
$ws.Range("K97").Value = 474.93332
$ws.Range("L97").Value = 568.375
$ws.Range("M97").Value = 21.06668000000002
$ws.Range("N97").Value = -1560.375
$ws.Range("H102").Value = 1210.0769
$ws.Range("I102").Value = 1217.6666
$ws.Range("J102").Value = 1193
$ws.Range("K102").Value = 1217.6666
$ws.Range("L102").Value = 1193
$ws.Range("M102").Value = 404.3334
$ws.Range("N102").Value = -4437
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1899.6923
$ws.Range("I7").Value = 1817.3636
$ws.Range("K7").Value = 1817.3636
$ws.Range("M7").Value = -1705.3636
$ws.Range("H40").Value = 7478.5713
$ws.Range("I40").Value = 2000
$ws.Range("K40").Value = 2000
$ws.Range("M40").Value = -1864
$ws.Range("H46").Value = 5750
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 6428.5713
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 6428.5713
$ws.Range("M46").Value = -812
$ws.Range("N46").Value = -6804.5713
$ws.Range("H93").Value = 788.36365
$ws.Range("I93").Value = 787.8
$ws.Range("J93").Value = 794
$ws.Range("K93").Value = 787.8
$ws.Range("L93").Value = 794
$ws.Range("M93").Value = 460.2
$ws.Range("N93").Value = -3290
$ws.Range("H100").Value = 1167.1666
$ws.Range("I100").Value = 1029.4286
$ws.Range("K100").Value = 1029.4286
$ws.Range("M100").Value = -488.4286
$ws.Range("H126").Value = 1899.6923
$ws.Range("I126").Value = 1817.3636
$ws.Range("K126").Value = 5452.0908
$ws.Range("M126").Value = -2982.0908

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 500
$ws.Range("I43").Value = 500
$ws.Range("K43").Value = 500
$ws.Range("M43").Value = -351
$ws.Range("H122").Value = 14707506
$ws.Range("I122").Value = 20835216
$ws.Range("J122").Value = 1001.8
$ws.Range("K122").Value = 62505648
$ws.Range("L122").Value = 3005.4
$ws.Range("M122").Value = -62503198
$ws.Range("N122").Value = -7905.4
$ws.Range("H125").Value = 49999.5
$ws.Range("J125").Value = 49999.5
$ws.Range("L125").Value = 49999.5
$ws.Range("N125").Value = -59839.5
$ws.Range("H136").Value = 836
$ws.Range("I136").Value = 737.6667
$ws.Range("K136").Value = 2213.0001
$ws.Range("M136").Value = 336.9998999999998

Write-Output "Applied all Kujata_Profits changes"